{"js": "// Seminar Notes.docx edit \u2014 adds the new seminar-assignment intro\n// paragraph (+ blank line) at the top of the document, and appends the\n// two new seminar write-ups (+ blank lines) at the end, matching the\n// target XML diff.\n\nconst body = context.document.body;\n\n// --- 1. Prepend the assignment-instructions paragraph + blank line ---\n\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst firstPara = paras.items[0];\n\nconst introText =\n  \"Attend at least five seminars/webinars related to the natural sciences \" +\n  \"(biology, chemistry, geology, physics, environment, etc). Check with \" +\n  \"the Ogden College departmental offices and/or websites for lists of \" +\n  \"upcoming seminars this semester. You must ask the speaker at least one \" +\n  \"question either during or after the seminar. Write a 1 paragraph \" +\n  \"synopsis of each seminar; include the title of the seminar, the \" +\n  \"speaker\\u2019s name and institution, time/date and location, the main \" +\n  \"point(s) of the talk, your question, and the speaker\\u2019s reply. It \" +\n  \"is highly recommended that you complete your synopsis on the same day \" +\n  \"as the seminar. You can also attend webinars from the American \" +\n  \"Chemical Society. Go to this website for a list of past and upcoming \" +\n  \"webinars. If you are not able to get your question in during the \" +\n  \"webinar, then try to contact the speaker via email, if possible\";\n\n// Insert the blank paragraph first, then the text paragraph before that\n// blank one, so the final order is: intro paragraph, blank paragraph,\n// original first paragraph (\"Q-Chem Webinar 65 ...\").\nconst blankBeforeFirst = firstPara.insertParagraph(\"\", \"Before\");\nblankBeforeFirst.insertParagraph(introText, \"Before\");\n\n// --- 2. Append the two new seminar write-ups (+ blank lines) at the end ---\n\nconst hhTitle =\n  \"Qualitative Investigation of the validity of an analytical Hodgkin-Huxley Model solution\";\nconst qmText = \"Quantum Mechanics on the cheap by Dr. Brenda Rubensteins\";\n\nconst hhPara = body.insertParagraph(hhTitle, \"End\");\nconst blankAfterHH = body.insertParagraph(\"\", \"End\");\nconst qmPara = body.insertParagraph(qmText, \"End\");\nconst finalBlank = body.insertParagraph(\"\", \"End\");\nawait context.sync();\n\n// Hanging indent: left indent 36pt (720 twips), first line -36pt (hanging 720 twips)\n// Set this AFTER the trailing blank paragraph already exists so the blank\n// paragraph does not inherit the indent formatting from qmPara.\nqmPara.leftIndent = 36;\nqmPara.firstLineIndent = -36;\n\nawait context.sync();\n", "ps1": "# Seminar Notes.docx edit \u2014 adds the new seminar-assignment intro\n# paragraph (+ blank line) at the top of the document, and appends the\n# two new seminar write-ups (+ blank lines) at the end, matching the\n# target XML diff.\n\n$d = $word.ActiveDocument\n\n# --- 1. Prepend the assignment-instructions paragraph + blank line ---\n\n$introText = \"Attend at least five seminars/webinars related to the natural sciences (biology, chemistry, geology, physics, environment, etc). Check with the Ogden College departmental offices and/or websites for lists of upcoming seminars this semester. You must ask the speaker at least one question either during or after the seminar. Write a 1 paragraph synopsis of each seminar; include the title of the seminar, the speaker\" + [char]8217 + \"s name and institution, time/date and location, the main point(s) of the talk, your question, and the speaker\" + [char]8217 + \"s reply. It is highly recommended that you complete your synopsis on the same day as the seminar. You can also attend webinars from the American Chemical Society. Go to this website for a list of past and upcoming webinars. If you are not able to get your question in during the webinar, then try to contact the speaker via email, if possible\"\n\n$firstPara = $d.Paragraphs(1)\n$firstPara.Range.InsertBefore($introText + [char]13 + [char]13)\n\n# --- 2. Append the two new seminar write-ups (+ blank lines) at the end ---\n\n$hhTitle = \"Qualitative Investigation of the validity of an analytical Hodgkin-Huxley Model solution\"\n$qmText = \"Quantum Mechanics on the cheap by Dr. Brenda Rubensteins\"\n\n$lastPara = $d.Paragraphs.Last\n$lastPara.Range.InsertAfter([char]13 + $hhTitle + [char]13 + [char]13 + $qmText + [char]13)\n\n# The newly-added \"Quantum Mechanics...\" paragraph is second-to-last\n# (the trailing blank paragraph is last). Give it the hanging indent:\n# left indent 36pt (720 twips), first line -36pt (hanging 720 twips).\n$n = $d.Paragraphs.Count\n$qmPara = $d.Paragraphs($n - 1)\n$qmPara.LeftIndent = 36\n$qmPara.FirstLineIndent = -36\n"}
